$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Add a default footer containing the standard built-in "Footer
# (Three Columns)" layout: a 3-column borderless table (left / center /
# right tab-style cells using the Header paragraph style) followed by an
# empty Footer-styled paragraph. Word materializes this via the footer
# story's Range.InsertXML (the only writable path into HeaderFooter
# content here; Range.WordOpenXML itself is read-only).
$footers = $sec.Footers
$ftr = $footers.Item(1)
$rng = $ftr.Range
$footerXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="TableNormal"/><w:bidiVisual w:val="0"/><w:tblW w:w="0" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3120"/><w:gridCol w:w="3120"/><w:gridCol w:w="3120"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:ind w:left="-115"/><w:jc w:val="left"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:ind w:right="-115"/><w:jc w:val="right"/></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="Footer"/><w:bidi w:val="0"/></w:pPr></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$rng.InsertXML($footerXml)

# --- Explicitly stamp the page as portrait, which makes Word persist
# w:orient="portrait" on <w:pgSz> instead of leaving it implicit.
$sec.PageSetup.Orientation = 0
